# Apply the changes described by the commit diff to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Students (Made by Luisa)" -> "Students"
$ws.Name = "Students"

# Update the view: drop the scrolled-to-row922 "topLeftCell" position and
# instead zoom the sheet to 130%, while keeping the existing selection
# (A941) intact.
$win = $excel.ActiveWindow
$win.Zoom = 130
